$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "298.83"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-2.28%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.65"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.71%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.117"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.54%"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.83%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.783"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.61%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.724"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "10.95%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.797"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.26%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9264"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.36%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1707"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.27%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07622"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.62%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07991"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.07%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03050"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.87%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09902"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.43%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001488"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.90%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04656"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2.44%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006150"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-4.78%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.457"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.53%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.51%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3291"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.50%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1335"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.59%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.566"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "8.25%"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-4.42%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001215"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.65%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004421"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.41%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001400"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "19.68%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001808"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "8.64%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01669"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "1.93%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04549"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.17%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006942"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-4.99%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1327"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-2.84%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002059"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-8.84%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01284"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-7.25%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006071"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.60%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7116"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-62.40%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-5.71%"

